# Update mods data [2025-12-07 15:08:30]
# Append a new daily data row (row 28) to the ModCounts sheet:
#   A28 = "2025/12/07" (text, like the other Date cells)
#   B28 = "逃离鸭科夫"   (text, same game name as every other row)
#   C28 = 1343          (numeric mod count)
# The new row should look/format just like the existing data rows
# (rows 3-27), which all share one cell style (centered alignment,
# general number format).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 28

# --- Column A: the date column is stored as literal text in this sheet
# (e.g. "2025/11/11"), not as a real Excel date. Typing a yyyy/mm/dd-style
# string straight into a General-formatted cell would get auto-parsed into
# a date serial, so force the cell to Text first, enter the value, then
# drop the now-unneeded "@" number format back off the cell (ClearFormats)
# so it ends up with the same plain/general formatting as its neighbours.
$ws.Range("A" + $row).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2025/12/07"
$ws.Range("A" + $row).ClearFormats()

# --- Column B: plain text game name, same as every other row.
$ws.Cells.Item($row, 2).Value = "逃离鸭科夫"

# --- Column C: numeric mod count.
$ws.Cells.Item($row, 3).Value = 1343

# --- Match the visual style of the existing data rows (centered
# horizontal/vertical alignment) for the whole new row.
$ws.Range("A" + $row + ":C" + $row).HorizontalAlignment = -4108
$ws.Range("A" + $row + ":C" + $row).VerticalAlignment = -4108
